$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function SetText($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

SetText $ws "D2" "63.777.79"
SetText $ws "E2" "  +0.81%  "
SetText $ws "D3" "3.143.24"
SetText $ws "E3" "  +1.02%  "
SetText $ws "E4" "  -0.07%  "
SetText $ws "D5" "587.33"
SetText $ws "E5" "  +0.37%  "
SetText $ws "D6" "145.48"
SetText $ws "E6" "  +0.37%  "
SetText $ws "E7" "  -0.01%  "
SetText $ws "D8" "3.137.97"
SetText $ws "E8" "  +1.17%  "
SetText $ws "D9" "0.530"
SetText $ws "E9" "  -0.07%  "
SetText $ws "E10" "  +7.15%  "
SetText $ws "E11" "  -0.09%  "
SetText $ws "E12" "  -2.10%  "
SetText $ws "D13" "0.0000247"
SetText $ws "E13" "  -0.19%  "
SetText $ws "D14" "36.93"
SetText $ws "E14" "  +3.81%  "
SetText $ws "D15" "3.664.13"
SetText $ws "E15" "  +1.05%  "
SetText $ws "E16" "  -1.41%  "
SetText $ws "B17" "WrappedBTC"
SetText $ws "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
SetText $ws "D17" "63.583.24"
SetText $ws "E17" "  +0.62%  "
SetText $ws "B18" "WrappedEther"
SetText $ws "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
SetText $ws "D18" "3.140.80"
SetText $ws "E18" "  +0.94%  "
SetText $ws "E19" "  -0.89%  "
SetText $ws "D20" "463.39"
SetText $ws "E20" "  -0.79%  "
SetText $ws "D21" "14.26"
SetText $ws "E21" "  +0.78%  "
SetText $ws "D22" "0.730"
SetText $ws "E22" "  +0.46%  "
SetText $ws "E23" "  -1.13%  "
SetText $ws "D24" "12.98"
SetText $ws "E24" "  -2.56%  "
SetText $ws "D25" "81.38"
SetText $ws "E25" "  -0.90%  "
SetText $ws "E26" "  +1.12%  "
SetText $ws "E27" "  +0.12%  "
SetText $ws "D28" "9.16"
SetText $ws "E28" "  +7.21%  "
SetText $ws "E29" "  +0.09%  "
SetText $ws "B30" "ImmutableX"
SetText $ws "C30" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
SetText $ws "D30" "2.21"
SetText $ws "E30" "  -0.37%  "
SetText $ws "B31" "FirstDigitalUSD"
SetText $ws "C31" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
SetText $ws "D31" "1.00"
SetText $ws "E31" "  -0.10%  "
SetText $ws "D32" "6.97"
SetText $ws "E32" "  +1.20%  "
SetText $ws "D33" "26.95"
SetText $ws "E33" "  +0.07%  "
SetText $ws "E34" "  +0.23%  "
SetText $ws "D35" "0.0₃0855"
SetText $ws "E35" "  -1.07%  "
SetText $ws "E36" "  -0.37%  "
SetText $ws "D37" "3.37"
SetText $ws "E37" "  +1.49%  "
SetText $ws "E38" "  -4.38%  "
SetText $ws "E39" "  -1.02%  "
SetText $ws "D40" "50.81"
SetText $ws "E40" "  +0.54%  "
SetText $ws "D41" "440.15"
SetText $ws "E41" "  +0.13%  "
SetText $ws "D42" "8.80"
SetText $ws "E42" "  +0.95%  "
SetText $ws "E43" "  +0.28%  "
SetText $ws "D44" "2.911.81"
SetText $ws "E44" "  -0.35%  "
SetText $ws "D45" "0.276"
SetText $ws "E45" "  -1.09%  "
SetText $ws "E46" "  -1.66%  "
SetText $ws "D47" "36.74"
SetText $ws "E47" "  +4.97%  "
SetText $ws "D48" "125.71"
SetText $ws "E48" "  +2.14%  "
SetText $ws "E50" "  -0.66%  "
SetText $ws "D51" "24.43"
SetText $ws "E51" "  -0.98%  "
